{"js": "const ooxml = \"<?xml version=\\\"1.0\\\" encoding=\\\"UTF-8\\\" standalone=\\\"yes\\\"?><pkg:package xmlns:pkg=\\\"http://schemas.microsoft.com/office/2006/xmlPackage\\\"><pkg:part pkg:name=\\\"/word/document.xml\\\" pkg:contentType=\\\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\\\"><pkg:xmlData><w:document xmlns:w=\\\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\\\"><w:body><w:p><w:r><w:t>Scoliosis</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\"> is a pathological curvature of the spine </w:t></w:r><w:r><w:t>affecting</w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\"> a few percent of the population. </w:t></w:r><w:r><w:t>The curvature tends to develop throughout growth. Therefore regular monitoring is important to make sure the disease</w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\"> doesn\\u2019t</w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\"> cause health problems.  For years, the gold-stan</w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\">dard approach to this has been </w:t></w:r><w:r><w:br/></w:r><w:r><w:t>X-raying the back, and measuring the Cobb angle, the greatest angle between the end-plates of any two vertebrae. As you might expect, regularly X-raying adolescents creates health risks. This has motivated research into spatially tracked ultrasound as a means to quantify the disease.</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\\\"preserve\\\">A number of preliminary studies have been done to assess the accuracy and validity of using tracked ultrasound. </w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\">They </w:t></w:r><w:r><w:t>consist of</w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\"> finding, and placing points at</w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\"> anatomic landmarks</w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\"> like the transverse processes, shown on my slide</w:t></w:r><w:r><w:t>.</w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\"> Curvatures are computed from these points\\u2019 </w:t></w:r><w:r><w:t>locations, and compared to the gold-standard</w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\"> X-ray</w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\"> measurement</w:t></w:r><w:r><w:t>. Regardless of whether these studies have been done on phantom models, or living patients, I suspect that their land</w:t></w:r><w:r><w:t>mark data is somewhat idealized, like my left picture.</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\\\"preserve\\\">I am only aware of </w:t></w:r><w:r><w:t>validation work done on mild to moderate cases of scoliosis, Cobb angles mostly up to 45</w:t></w:r><w:r><w:rPr><w:vertAlign w:val=\\\"superscript\\\"/></w:rPr><w:t>o</w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\">. There are also idealizations associated with using phantom models, and wide-range transducer setups. Some of these authors themselves acknowledge being unable to locate anatomic landmarks in </w:t></w:r><w:r><w:t>a few</w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\"> </w:t></w:r><w:r><w:t>places</w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\">. </w:t></w:r><w:r><w:t>Apparently, a quantification method for working with a scope including imperfect data, and cases worse than 45</w:t></w:r><w:r><w:rPr><w:vertAlign w:val=\\\"superscript\\\"/></w:rPr><w:t>o</w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\"> or with other complications is </w:t></w:r><w:r><w:t>needed</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:r><w:t>I claim that a method using neural networks can be developed addressing this need.</w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\"> Neural networks are known for their robustness and accuracy, indispensable virtues in clinical settings.</w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\"> I have 124 sets of scoliotic patients\\u2019 transverse processes </w:t></w:r><w:r><w:t>locations</w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\"> from CT-derived models</w:t></w:r><w:r><w:t>, ranging from mild to severe</w:t></w:r><w:r><w:t>. The accuracy of CT makes this a natural ground-truth from which to extract the correct curvature. I will then programmatically degrade the data, introduci</w:t></w:r><w:r><w:t>ng errors expected in ultrasound</w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\">: noise, missing points, and misplaced points. A pre-processing module will be </w:t></w:r><w:r><w:t>developed</w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\"> to repair the data sets so they can be used to train a neural network.</w:t></w:r></w:p><w:p><w:r><w:t>After being trained on the basis of the difference between the network\\u2019s output and the correct angle, the network can be similarly evaluated.</w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\"> The average curvature estimate error in degrees, and other statistics, will be collected for various experimental setups.</w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\"> </w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\">I will vary the amounts of input data error, use different training set sizes, networks architectures and functionalities, </w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\">and so on. </w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\">With the results of these </w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\">planned </w:t></w:r><w:r><w:t>experiments, I hope to demonstrate that the method produces curvature estimates comparable to those of</w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\"> </w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\">current studies, and within clinically </w:t></w:r><w:r><w:t>acceptable limits of error, including for these difficult cases</w:t></w:r><w:bookmarkStart w:id=\\\"0\\\" w:name=\\\"_GoBack\\\"/><w:bookmarkEnd w:id=\\\"0\\\"/><w:r><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\";\n\nconst body = context.document.body;\nbody.insertOoxml(ooxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$xml = @'\n<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:r><w:t>Scoliosis</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space=\"preserve\"> is a pathological curvature of the spine </w:t></w:r><w:r><w:t>affecting</w:t></w:r><w:r><w:t xml:space=\"preserve\"> a few percent of the population. </w:t></w:r><w:r><w:t>The curvature tends to develop throughout growth. Therefore regular monitoring is important to make sure the disease</w:t></w:r><w:r><w:t xml:space=\"preserve\"> doesn\u2019t</w:t></w:r><w:r><w:t xml:space=\"preserve\"> cause health problems.  For years, the gold-stan</w:t></w:r><w:r><w:t xml:space=\"preserve\">dard approach to this has been </w:t></w:r><w:r><w:br/></w:r><w:r><w:t>X-raying the back, and measuring the Cobb angle, the greatest angle between the end-plates of any two vertebrae. As you might expect, regularly X-raying adolescents creates health risks. This has motivated research into spatially tracked ultrasound as a means to quantify the disease.</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">A number of preliminary studies have been done to assess the accuracy and validity of using tracked ultrasound. </w:t></w:r><w:r><w:t xml:space=\"preserve\">They </w:t></w:r><w:r><w:t>consist of</w:t></w:r><w:r><w:t xml:space=\"preserve\"> finding, and placing points at</w:t></w:r><w:r><w:t xml:space=\"preserve\"> anatomic landmarks</w:t></w:r><w:r><w:t xml:space=\"preserve\"> like the transverse processes, shown on my slide</w:t></w:r><w:r><w:t>.</w:t></w:r><w:r><w:t xml:space=\"preserve\"> Curvatures are computed from these points\u2019 </w:t></w:r><w:r><w:t>locations, and compared to the gold-standard</w:t></w:r><w:r><w:t xml:space=\"preserve\"> X-ray</w:t></w:r><w:r><w:t xml:space=\"preserve\"> measurement</w:t></w:r><w:r><w:t>. Regardless of whether these studies have been done on phantom models, or living patients, I suspect that their land</w:t></w:r><w:r><w:t>mark data is somewhat idealized, like my left picture.</w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">I am only aware of </w:t></w:r><w:r><w:t>validation work done on mild to moderate cases of scoliosis, Cobb angles mostly up to 45</w:t></w:r><w:r><w:rPr><w:vertAlign w:val=\"superscript\"/></w:rPr><w:t>o</w:t></w:r><w:r><w:t xml:space=\"preserve\">. There are also idealizations associated with using phantom models, and wide-range transducer setups. Some of these authors themselves acknowledge being unable to locate anatomic landmarks in </w:t></w:r><w:r><w:t>a few</w:t></w:r><w:r><w:t xml:space=\"preserve\"> </w:t></w:r><w:r><w:t>places</w:t></w:r><w:r><w:t xml:space=\"preserve\">. </w:t></w:r><w:r><w:t>Apparently, a quantification method for working with a scope including imperfect data, and cases worse than 45</w:t></w:r><w:r><w:rPr><w:vertAlign w:val=\"superscript\"/></w:rPr><w:t>o</w:t></w:r><w:r><w:t xml:space=\"preserve\"> or with other complications is </w:t></w:r><w:r><w:t>needed</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:r><w:t>I claim that a method using neural networks can be developed addressing this need.</w:t></w:r><w:r><w:t xml:space=\"preserve\"> Neural networks are known for their robustness and accuracy, indispensable virtues in clinical settings.</w:t></w:r><w:r><w:t xml:space=\"preserve\"> I have 124 sets of scoliotic patients\u2019 transverse processes </w:t></w:r><w:r><w:t>locations</w:t></w:r><w:r><w:t xml:space=\"preserve\"> from CT-derived models</w:t></w:r><w:r><w:t>, ranging from mild to severe</w:t></w:r><w:r><w:t>. The accuracy of CT makes this a natural ground-truth from which to extract the correct curvature. I will then programmatically degrade the data, introduci</w:t></w:r><w:r><w:t>ng errors expected in ultrasound</w:t></w:r><w:r><w:t xml:space=\"preserve\">: noise, missing points, and misplaced points. A pre-processing module will be </w:t></w:r><w:r><w:t>developed</w:t></w:r><w:r><w:t xml:space=\"preserve\"> to repair the data sets so they can be used to train a neural network.</w:t></w:r></w:p><w:p><w:r><w:t>After being trained on the basis of the difference between the network\u2019s output and the correct angle, the network can be similarly evaluated.</w:t></w:r><w:r><w:t xml:space=\"preserve\"> The average curvature estimate error in degrees, and other statistics, will be collected for various experimental setups.</w:t></w:r><w:r><w:t xml:space=\"preserve\"> </w:t></w:r><w:r><w:t xml:space=\"preserve\">I will vary the amounts of input data error, use different training set sizes, networks architectures and functionalities, </w:t></w:r><w:r><w:t xml:space=\"preserve\">and so on. </w:t></w:r><w:r><w:t xml:space=\"preserve\">With the results of these </w:t></w:r><w:r><w:t xml:space=\"preserve\">planned </w:t></w:r><w:r><w:t>experiments, I hope to demonstrate that the method produces curvature estimates comparable to those of</w:t></w:r><w:r><w:t xml:space=\"preserve\"> </w:t></w:r><w:r><w:t xml:space=\"preserve\">current studies, and within clinically </w:t></w:r><w:r><w:t>acceptable limits of error, including for these difficult cases</w:t></w:r><w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/><w:r><w:t>.</w:t></w:r></w:p>\n'@\n\n$d.Content.InsertXML($xml)\n"}
